# dialbb no-code scenario template: casing fix for the ChatGPT NER
# placeholder (#NE_PERSON -> #NE_person) plus the row-height/selection
# touch-ups that came along with the re-save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content fix: NE_PERSON -> NE_person (casing) -------------------------
# F2 holds the condition string; C4 holds the system-utterance template
# that references the same named-entity placeholder.
$ws.Range("F2").Value = '#NE_person!=""'
$ws.Range("C4").Value = 'Thank you {#NE_person}! Let me ask you about sandwich. Do you have sandwiches very often?'

# --- row heights (re-wrap heights after the re-save) -----------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 60

# --- selection moved from G2 to C5 -----------------------------------------
$ws.Range("C5").Select()
